# Auto-generated edit script applying the Alpha_Profits.xlsx diff
# Updates/clears/adds numeric cells across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the target revision described in the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J33").Value = 727.75
$ws.Range("K33").Value = 33659.066
$ws.Range("M33").Value = -33430.066
$ws.Range("L33").Value = 727.75
$ws.Range("N33").Value = -1185.75
$ws.Range("I33").Value = 33659.066
$ws.Range("H33").Value = 26726.158
$ws.Range("J51").Value = 8326.333000000001
$ws.Range("L51").Value = 8326.333000000001
$ws.Range("H51").Value = 25439.777
$ws.Range("N51").Value = -9294.333000000001
$ws.Range("N62").Value = -4218
$ws.Range("L62").Value = 2970
$ws.Range("H62").Value = 3024.8333
$ws.Range("J62").Value = 2970
$ws.Range("N65").Value = -21090
$ws.Range("H65").Value = 3024.8333
$ws.Range("J65").Value = 2970
$ws.Range("L65").Value = 14850
$ws.Range("H99").Value = 946.3077
$ws.Range("M99").Value = 461.09095
$ws.Range("I99").Value = 345.63635
$ws.Range("J99").Value = 4250
$ws.Range("K99").Value = 1036.90905
$ws.Range("L99").Value = 12750
$ws.Range("N99").Value = -15746
$ws.Range("N101").Value = -12364
$ws.Range("L101").Value = 9120
$ws.Range("J101").Value = 3040
$ws.Range("H101").Value = 1483.3
$ws.Range("J110").Value = 45695
$ws.Range("N110").Value = -53875
$ws.Range("H110").Value = 45695
$ws.Range("L110").Value = 45695
$ws.Range("J111").Value = 1420.3334
$ws.Range("N111").Value = -10395.0002
$ws.Range("M111").Value = -255.0001999999999
$ws.Range("H111").Value = 1211.6666
$ws.Range("L111").Value = 4261.0002
$ws.Range("I111").Value = 1107.3334
$ws.Range("K111").Value = 3322.0002
$ws.Range("I114").Value = 65000
$ws.Range("H114").Value = 66750
$ws.Range("M114").Value = -60661
$ws.Range("K114").Value = 65000
$ws.Range("I132").Value = 94538.27
$ws.Range("K132").Value = 283614.81
$ws.Range("H132").Value = 75287.14
$ws.Range("M132").Value = -281084.81

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 9271.733
$ws.Range("K2").Value = 9271.733
$ws.Range("H2").Value = 8305.65
$ws.Range("M2").Value = -9158.733
$ws.Range("M45").Value = -1558
$ws.Range("L45").Value = 2054.5
$ws.Range("J45").Value = 2054.5
$ws.Range("H45").Value = 2014.6666
$ws.Range("I45").Value = 1935
$ws.Range("N45").Value = -2808.5
$ws.Range("K45").Value = 1935
$ws.Range("K74").Value = 1684379.6
$ws.Range("M74").Value = -1683505.6
$ws.Range("I74").Value = 1684379.6
$ws.Range("H74").Value = 2572869.2
$ws.Range("I77").Value = 1684379.6
$ws.Range("H77").Value = 2572869.2
$ws.Range("K77").Value = 8421898
$ws.Range("M77").Value = -8417530
$ws.Range("J92").Value = 30049
$ws.Range("L92").Value = 30049
$ws.Range("N92").Value = -35041
$ws.Range("H92").Value = 30049
$ws.Range("M102").Value = -6.285800000000108
$ws.Range("H102").Value = 1673
$ws.Range("K102").Value = 1628.2858
$ws.Range("I102").Value = 1628.2858
$ws.Range("I110").Value = 8991.362999999999
$ws.Range("K110").Value = 8991.362999999999
$ws.Range("M110").Value = -6946.362999999999
$ws.Range("H110").Value = 9424.357
$ws.Range("H116").Value = 8305.65
$ws.Range("I116").Value = 9271.733
$ws.Range("M116").Value = -6977.733
$ws.Range("K116").Value = 9271.733

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 9271.733
$ws.Range("M3").Value = -9157.733
$ws.Range("H3").Value = 8305.65
$ws.Range("I3").Value = 9271.733
$ws.Range("L93").Value = 44983
$ws.Range("J93").Value = 44983
$ws.Range("N93").Value = -48727
$ws.Range("H93").Value = 44983
$ws.Range("K94").Value = 25181.875
$ws.Range("M94").Value = -24730.875
$ws.Range("H94").Value = 18431.182
$ws.Range("I94").Value = 25181.875
$ws.Range("H99").Value = 4061.4
$ws.Range("M99").Value = -2563.4
$ws.Range("I99").Value = 4061.4
$ws.Range("K99").Value = 4061.4
$ws.Range("I107").Value = 3040.8
$ws.Range("J107").Value = 85000
$ws.Range("L107").Value = 85000
$ws.Range("K107").Value = 3040.8
$ws.Range("H107").Value = 33775.5
$ws.Range("N107").Value = -88840
$ws.Range("M107").Value = -1120.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L18").Value = 26930.5
$ws.Range("H18").Value = 26930.5
$ws.Range("N18").Value = -27390.5
$ws.Range("J18").Value = 26930.5
$ws.Range("L94").Value = 767.1111
$ws.Range("N94").Value = -1669.1111
$ws.Range("H94").Value = 808.4545000000001
$ws.Range("J94").Value = 767.1111
$ws.Range("I107").Value = 4909
$ws.Range("J107").Value = 4310.6665
$ws.Range("L107").Value = 4310.6665
$ws.Range("K107").Value = 4909
$ws.Range("H107").Value = 4582.636
$ws.Range("N107").Value = -8150.6665
$ws.Range("M107").Value = -2989
$ws.Range("I122").Value = 100010
$ws.Range("M122").Value = -297580
$ws.Range("K122").Value = 300030
$ws.Range("H122").Value = 175005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H31").Value = 600
$ws.Range("H121").Value = 1751947.6
$ws.Range("J121").Value = 3335668.5
$ws.Range("N121").Value = -10009625.5
$ws.Range("I121").Value = 168226.83
$ws.Range("M121").Value = -503370.49
$ws.Range("L121").Value = 10007005.5
$ws.Range("K121").Value = 504680.49

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N70").Value = -8706.166499999999
$ws.Range("J70").Value = 8166.1665
$ws.Range("L70").Value = 8166.1665
$ws.Range("H70").Value = 6621.28
$ws.Range("K70").Value = 6133.421
$ws.Range("M70").Value = -5863.421
$ws.Range("I70").Value = 6133.421
$ws.Range("I73").Value = 6133.421
$ws.Range("K73").Value = 6133.421
$ws.Range("L73").Value = 8166.1665
$ws.Range("H73").Value = 6621.28
$ws.Range("M73").Value = -5197.421
$ws.Range("N73").Value = -10038.1665
$ws.Range("J73").Value = 8166.1665
$ws.Range("M102").ClearContents()
$ws.Range("H102").Value = 4499.5
$ws.Range("K102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("I122").Value = 2631.158
$ws.Range("M122").Value = -5443.474
$ws.Range("K122").Value = 7893.474
$ws.Range("H122").Value = 2290.5173
$ws.Range("N126").Value = -15440
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -50327.75
$ws.Range("J126").Value = 3500
$ws.Range("I126").Value = 17599.25
$ws.Range("K126").Value = 52797.75
$ws.Range("H126").Value = 11556.714

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K61").Value = 2178.2856
$ws.Range("M61").Value = -1976.2856
$ws.Range("H61").Value = 3471.889
$ws.Range("I61").Value = 2178.2856
$ws.Range("H96").Value = 49153.5
$ws.Range("N96").Value = -54645.5
$ws.Range("L96").Value = 49153.5
$ws.Range("J96").Value = 49153.5
$ws.Range("H113").Value = 3471.889
$ws.Range("M113").Value = -8.285600000000159
$ws.Range("I113").Value = 2178.2856
$ws.Range("K113").Value = 2178.2856
$ws.Range("I122").Value = 2857.375
$ws.Range("M122").Value = -6122.125
$ws.Range("K122").Value = 8572.125
$ws.Range("H122").Value = 3222.8635
$ws.Range("I132").Value = 2775.9375
$ws.Range("K132").Value = 8327.8125
$ws.Range("H132").Value = 3205.4211
$ws.Range("M132").Value = -5797.8125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I107").Value = 1002.3333
$ws.Range("K107").Value = 3006.9999
$ws.Range("H107").Value = 2268.2
$ws.Range("M107").Value = -1086.9999
$ws.Range("J113").Value = 8502.4
$ws.Range("H113").Value = 6313.857
$ws.Range("M113").Value = -357.5
$ws.Range("N113").Value = -29847.2
$ws.Range("I113").Value = 842.5
$ws.Range("K113").Value = 2527.5
$ws.Range("L113").Value = 25507.2
$ws.Range("I122").Value = 5167.125
$ws.Range("M122").Value = -13051.375
$ws.Range("K122").Value = 15501.375
$ws.Range("H122").Value = 4142.647
